$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 value changes from 11.023.210/0001-11 to 27.988.301/0001-21
$ws.Range("A2").Value = "27.988.301/0001-21"

# Rows 11 and 12 get cleared (were 28.860.597/0001-63 and 42.856.348/0001-25)
$ws.Range("A11").ClearContents()
$ws.Range("A12").ClearContents()

# Update selection to A2:A10, active cell A2
$ws.Range("A2:A10").Select()
